# Add two new inventory rows (2 and 3) to the "Inventory" sheet, matching
# the existing header row's 15 columns (A:O). This grows the sheet's
# dimension from A1:O1 to A1:O3.
#
# Columns D, F, G, H, M are genuine numbers; the remaining columns
# (A, B, C, E, I, J, K, L, N, O) are text. Some of the text values look
# like numbers ("4.73", "15.23", "0.71", "15.93", "0.00") so we must force
# them to stay text instead of being auto-coerced to numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")
$ws.Activate()

$numericTextCols = 9, 10, 11, 12, 14, 15
$textValues = @{
    9  = "4.73"
    10 = "15.23"
    11 = "0.71"
    12 = "15.93"
    14 = "0.00"
    15 = "0.00"
}

foreach ($r in 2, 3) {
    $ws.Cells.Item($r, 1).Value = "glop"
    $ws.Cells.Item($r, 2).Value = "stock"
    $ws.Cells.Item($r, 3).Value = "eth"
    $ws.Cells.Item($r, 4).Value = 12
    $ws.Cells.Item($r, 5).Value = "pcs"
    $ws.Cells.Item($r, 6).Value = 23
    $ws.Cells.Item($r, 7).Value = 921321546465
    $ws.Cells.Item($r, 8).Value = 10.5

    # Force these number-looking values to remain plain text cells.
    foreach ($c in $numericTextCols) {
        $ws.Cells.Item($r, $c).NumberFormat = "@"
        $ws.Cells.Item($r, $c).Value = $textValues[$c]
    }

    $ws.Cells.Item($r, 13).Value = 23
}

# The "@" text format left a non-default style on the text cells above;
# paste the header cell's (default) formatting back over them so the new
# rows don't carry a stray explicit number format.
$ws.Cells.Item(1, 1).Copy()
foreach ($r in 2, 3) {
    foreach ($c in $numericTextCols) {
        $ws.Cells.Item($r, $c).PasteSpecial(-4122)  # xlPasteFormats
    }
}
$excel.CutCopyMode = $false
